# Scheduled market-data refresh: recompute currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) for the affected leve rows across all crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 1004.2308
$ws.Range("I33").Value = 1037.9584
$ws.Range("J33").Value = 599.5
$ws.Range("K33").Value = 1037.9584
$ws.Range("L33").Value = 599.5
$ws.Range("M33").Value = -808.9584
$ws.Range("N33").Value = -1057.5

# Row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value = 11113432
$ws.Range("I51").Value = 3771
$ws.Range("J51").Value = 15874716
$ws.Range("K51").Value = 3771
$ws.Range("L51").Value = 15874716
$ws.Range("M51").Value = -3287
$ws.Range("N51").Value = -15875684

# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 13995
$ws.Range("I88").Value = 21500
$ws.Range("J88").Value = 6490
$ws.Range("K88").Value = 21500
$ws.Range("L88").Value = 6490
$ws.Range("M88").Value = -21094
$ws.Range("N88").Value = -7302

# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 13995
$ws.Range("I91").Value = 21500
$ws.Range("J91").Value = 6490
$ws.Range("K91").Value = 21500
$ws.Range("L91").Value = 6490
$ws.Range("M91").Value = -20096
$ws.Range("N91").Value = -9298

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 2693.1304
$ws.Range("I113").Value = 2266.6667
$ws.Range("J113").Value = 2843.647
$ws.Range("K113").Value = 2266.6667
$ws.Range("L113").Value = 2843.647
$ws.Range("M113").Value = 987.3332999999998
$ws.Range("N113").Value = -9351.647000000001

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 8588144
$ws.Range("I116").Value = 10121112
$ws.Range("K116").Value = 10121112
$ws.Range("M116").Value = -10117670

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 5158.9697
$ws.Range("I132").Value = 5903.2354
$ws.Range("J132").Value = 4368.1875
$ws.Range("K132").Value = 17709.7062
$ws.Range("L132").Value = 13104.5625
$ws.Range("M132").Value = -15179.7062
$ws.Range("N132").Value = -18164.5625

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 1088.0769
$ws.Range("I61").Value = 941.375
$ws.Range("J61").Value = 1322.8
$ws.Range("K61").Value = 941.375
$ws.Range("L61").Value = 1322.8
$ws.Range("M61").Value = -729.375
$ws.Range("N61").Value = -1746.8

# Row 98: Greaving / Doman Iron Greaves of Maiming
$ws.Range("H98").Value = 46903.332
$ws.Range("J98").Value = 46903.332
$ws.Range("L98").Value = 46903.332
$ws.Range("N98").Value = -52893.332

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2714
$ws.Range("I122").Value = 1949
$ws.Range("K122").Value = 5847
$ws.Range("M122").Value = -3397

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1596.16
$ws.Range("I132").Value = 1258.1945
$ws.Range("J132").Value = 2465.2144
$ws.Range("K132").Value = 3774.5835
$ws.Range("L132").Value = 7395.6432
$ws.Range("M132").Value = -1244.5835
$ws.Range("N132").Value = -12455.6432

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1088.0769
$ws.Range("I136").Value = 941.375
$ws.Range("J136").Value = 1322.8
$ws.Range("K136").Value = 2824.125
$ws.Range("L136").Value = 3968.4
$ws.Range("M136").Value = -274.125
$ws.Range("N136").Value = -9068.4

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 20.5
$ws.Range("I22").Value = 20.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 20.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 152.5
$ws.Range("N22").ClearContents()

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1660
$ws.Range("I99").Value = 1666.6666
$ws.Range("J99").Value = 1650
$ws.Range("K99").Value = 1666.6666
$ws.Range("L99").Value = 1650
$ws.Range("M99").Value = -168.6666
$ws.Range("N99").Value = -4646

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1112.9807
$ws.Range("I134").Value = 1033.2565
$ws.Range("J134").Value = 1352.1538
$ws.Range("K134").Value = 3099.7695
$ws.Range("L134").Value = 4056.4614
$ws.Range("M134").Value = -564.7694999999999
$ws.Range("N134").Value = -9126.4614

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 1380.9062
$ws.Range("I58").Value = 759.8461
$ws.Range("J58").Value = 1805.8422
$ws.Range("K58").Value = 759.8461
$ws.Range("L58").Value = 1805.8422
$ws.Range("M58").Value = -556.8461
$ws.Range("N58").Value = -2211.8422

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 2103505
$ws.Range("I99").Value = 3970569.2
$ws.Range("J99").Value = 3057.5
$ws.Range("K99").Value = 3970569.2
$ws.Range("L99").Value = 3057.5
$ws.Range("M99").Value = -3969071.2
$ws.Range("N99").Value = -6053.5

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 2103505
$ws.Range("I126").Value = 3970569.2
$ws.Range("J126").Value = 3057.5
$ws.Range("K126").Value = 11911707.6
$ws.Range("L126").Value = 9172.5
$ws.Range("M126").Value = -11909237.6
$ws.Range("N126").Value = -14112.5

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 1213.76
$ws.Range("I132").Value = 1048.9231
$ws.Range("J132").Value = 1392.3334
$ws.Range("K132").Value = 3146.7693
$ws.Range("L132").Value = 4177.0002
$ws.Range("M132").Value = -616.7692999999999
$ws.Range("N132").Value = -9237.0002

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1583.1305
$ws.Range("I134").Value = 1278.9333
$ws.Range("K134").Value = 3836.7999
$ws.Range("M134").Value = -1301.7999

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 1380.9062
$ws.Range("I136").Value = 759.8461
$ws.Range("J136").Value = 1805.8422
$ws.Range("K136").Value = 2279.5383
$ws.Range("L136").Value = 5417.5266
$ws.Range("M136").Value = 270.4616999999998
$ws.Range("N136").Value = -10517.5266

$ws = $wb.Worksheets.Item("CUL")
# Row 98: Sweet Kiss of Death / Rice Vinegar
$ws.Range("H98").Value = 511.92307
$ws.Range("J98").Value = 551.25
$ws.Range("L98").Value = 1653.75
$ws.Range("N98").Value = -4649.75

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 59366.855
$ws.Range("I80").Value = 224264.44
$ws.Range("J80").Value = 2286.923
$ws.Range("K80").Value = 224264.44
$ws.Range("L80").Value = 2286.923
$ws.Range("M80").Value = -223266.44
$ws.Range("N80").Value = -4282.923

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 59366.855
$ws.Range("I83").Value = 224264.44
$ws.Range("J83").Value = 2286.923
$ws.Range("K83").Value = 1121322.2
$ws.Range("L83").Value = 11434.615
$ws.Range("M83").Value = -1116330.2
$ws.Range("N83").Value = -21418.615

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 2900.2856
$ws.Range("I113").Value = 1072.25
$ws.Range("K113").Value = 1072.25
$ws.Range("M113").Value = 1097.75

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 6629.643
$ws.Range("I122").Value = 5466.6665
$ws.Range("J122").Value = 7501.875
$ws.Range("K122").Value = 16399.9995
$ws.Range("L122").Value = 22505.625
$ws.Range("M122").Value = -13949.9995
$ws.Range("N122").Value = -27405.625

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 1134.8788
$ws.Range("I132").Value = 898.64703
$ws.Range("J132").Value = 1385.875
$ws.Range("K132").Value = 2695.94109
$ws.Range("L132").Value = 4157.625
$ws.Range("M132").Value = -165.9410899999998
$ws.Range("N132").Value = -9217.625

$ws = $wb.Worksheets.Item("LTW")
# Row 95: Weathering Heights / Gyuki Leather Highboots of Striking
$ws.Range("H95").Value = 29000
$ws.Range("J95").Value = 29000
$ws.Range("L95").Value = 29000
$ws.Range("N95").Value = -34492

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 5508.3335
$ws.Range("I122").Value = 5554.5454
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 16663.6362
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -14213.6362
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 5692.8667
$ws.Range("I62").Value = 4242.857
$ws.Range("J62").Value = 6961.625
$ws.Range("K62").Value = 4242.857
$ws.Range("L62").Value = 6961.625
$ws.Range("M62").Value = -3618.857
$ws.Range("N62").Value = -8209.625

# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 5692.8667
$ws.Range("I65").Value = 4242.857
$ws.Range("J65").Value = 6961.625
$ws.Range("K65").Value = 21214.285
$ws.Range("L65").Value = 34808.125
$ws.Range("M65").Value = -18094.285
$ws.Range("N65").Value = -41048.125

# Row 97: Getting a Leg Up / Ruby Cotton Gaskins of Striking
$ws.Range("H97").Value = 32000
$ws.Range("J97").Value = 32000
$ws.Range("L97").Value = 32000
$ws.Range("N97").Value = -33982

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1710.1052
$ws.Range("I122").Value = 1320.4
$ws.Range("J122").Value = 2143.111
$ws.Range("K122").Value = 3961.2
$ws.Range("L122").Value = 6429.333
$ws.Range("M122").Value = -1511.2
$ws.Range("N122").Value = -11329.333

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1996.4286
$ws.Range("I132").Value = 1336.6666
$ws.Range("J132").Value = 2656.1904
$ws.Range("K132").Value = 4009.9998
$ws.Range("L132").Value = 7968.5712
$ws.Range("M132").Value = -1479.9998
$ws.Range("N132").Value = -13028.5712
